$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Update evidence/data values in rows 9-13 ---
# Row 9
$ws.Range("B9").Value = "105952330"
$ws.Range("C9").Value = "3052755808"
$ws.Range("D9").Value = "732111324709812"
$ws.Range("H9").Value = "850396952"

# Row 10
$ws.Range("B10").Value = "841794534"
$ws.Range("C10").Value = "3052755811"
$ws.Range("D10").Value = "732111324709813"
$ws.Range("H10").Value = "29389214"

# Row 11
$ws.Range("B11").Value = "543372417"
$ws.Range("C11").Value = "3052755808"
$ws.Range("D11").Value = "732111324709812"

# Row 12
$ws.Range("B12").Value = "922697306"
$ws.Range("C12").Value = "3052755812"
$ws.Range("D12").Value = "732111324709814"

# Row 13
$ws.Range("B13").Value = "121541180"
$ws.Range("C13").Value = "3052755815"
$ws.Range("D13").Value = "732111324709816"

# --- Remove rows 14:16 (evidence rows no longer captured) ---
$ws.Range("A14:K16").EntireRow.Delete()

# --- Bold the secondary-header style cells (new cellXfs entry) ---
$ws.Range("E4:G4").Font.Bold = $true
$ws.Range("C8:I8").Font.Bold = $true
$ws.Range("E12").Font.Bold = $true

# --- Adjust view: clear the scrolled topLeftCell and move selection ---
$ws.Range("B16").Select()
